$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# XLSForm cleanup: the deprecated "hidden" and "string" question types are
# replaced with "text" (optionally combined with appearance="hidden").

# Row 4: source (hidden -> text, appearance=hidden)
$ws.Range("A4").Value = "text"
$ws.Range("F4").Value = "hidden"

# Row 5: source_id (hidden -> text, appearance=hidden)
$ws.Range("A5").Value = "text"
$ws.Range("F5").Value = "hidden"

# Row 7: _id (string -> text)
$ws.Range("A7").Value = "text"

# Row 8: name (hidden -> text, appearance=hidden)
$ws.Range("A8").Value = "text"
$ws.Range("F8").Value = "hidden"

# Move the cursor/selection on the frozen-pane bottom-right view to B16.
$ws.Range("B16").Select()
